$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0, bold font + thin box border + centered/top aligned
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160

# A2 = 0, same formatting as B1 - copy the format over so the
# engine reuses the same cell style record instead of minting a new one
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0

# B2 = "disconnected_elements" plain text, default style
$ws.Range("B2").Value = "disconnected_elements"
